# Adicionados balancos concatenados em uma unica planilha.
# Extends the VAMO3 sheet with three new quarterly columns (V=31/12/2023,
# W=31/03/2024, X=30/06/2024), mirroring the existing layout in columns B:U.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (V1:X1): quarter-end dates, matching style of row 1 ---
$headerArr = New-Object "object[,]" 1,3
$headerArr[0,0] = "31/12/2023"
$headerArr[0,1] = "31/03/2024"
$headerArr[0,2] = "30/06/2024"
$ws.Range("V1:X1").Value = $headerArr
$ws.Range("U1").Copy()
$ws.Range("V1:X1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Numeric data rows 2:56 (V:X) ---
$data_2 = New-Object "object[,]" 55,3
$data_2[0,0] = 20808826.88
$data_2[0,1] = 22334959.616
$data_2[0,2] = 22657783.808
$data_2[1,0] = 5957884.928
$data_2[1,1] = 6218613.76
$data_2[1,2] = 5992064
$data_2[2,0] = 97768
$data_2[2,1] = 265779.008
$data_2[2,2] = 240148.992
$data_2[3,0] = 2196243.968
$data_2[3,1] = 2182206.976
$data_2[3,2] = 1847373.056
$data_2[4,0] = 982814.0159999999
$data_2[4,1] = 1081400.96
$data_2[4,2] = 932200
$data_2[5,0] = 1650612.992
$data_2[5,1] = 1545852.032
$data_2[5,2] = 1572849.024
$data_2[6,0] = 0
$data_2[6,1] = 0
$data_2[6,2] = 0
$data_2[7,0] = 479008
$data_2[7,1] = 472276.992
$data_2[7,2] = 490860.992
$data_2[8,0] = 18015
$data_2[8,1] = 90181
$data_2[8,2] = 87941
$data_2[9,0] = 533423.008
$data_2[9,1] = 580916.992
$data_2[9,2] = 820691.008
$data_2[10,0] = 963081.024
$data_2[10,1] = 920627.008
$data_2[10,2] = 821881.024
$data_2[11,0] = 10950
$data_2[11,1] = 10899
$data_2[11,2] = 10638
$data_2[12,0] = 0
$data_2[12,1] = 0
$data_2[12,2] = 0
$data_2[13,0] = 0
$data_2[13,1] = 0
$data_2[13,2] = 0
$data_2[14,0] = 55511
$data_2[14,1] = 76841
$data_2[14,2] = 85961
$data_2[15,0] = 0
$data_2[15,1] = 0
$data_2[15,2] = 0
$data_2[16,0] = 0
$data_2[16,1] = 0
$data_2[16,2] = 0
$data_2[17,0] = 177600
$data_2[17,1] = 194224.992
$data_2[17,2] = 216142
$data_2[18,0] = 0
$data_2[18,1] = 0
$data_2[18,2] = 0
$data_2[19,0] = 0
$data_2[19,1] = 0
$data_2[19,2] = 0
$data_2[20,0] = 0
$data_2[20,1] = 0
$data_2[20,2] = 0
$data_2[21,0] = 13381557.248
$data_2[21,1] = 14694810.624
$data_2[21,2] = 15350112.256
$data_2[22,0] = 506303.008
$data_2[22,1] = 500908.992
$data_2[22,2] = 493727.008
$data_2[23,0] = 0
$data_2[23,1] = 0
$data_2[23,2] = 0
$data_2[24,0] = 20808826.88
$data_2[24,1] = 22334959.616
$data_2[24,2] = 22657783.808
$data_2[25,0] = 3412500.992
$data_2[25,1] = 4636029.952
$data_2[25,2] = 4558027.776
$data_2[26,0] = 72819
$data_2[26,1] = 61358
$data_2[26,2] = 76871
$data_2[27,0] = 1214952.96
$data_2[27,1] = 1638967.04
$data_2[27,2] = 1832162.944
$data_2[28,0] = 43224
$data_2[28,1] = 45024
$data_2[28,2] = 46405
$data_2[29,0] = 881625.024
$data_2[29,1] = 1572930.048
$data_2[29,2] = 1673122.048
$data_2[30,0] = 0
$data_2[30,1] = 0
$data_2[30,2] = 0
$data_2[31,0] = 300174.016
$data_2[31,1] = 302000
$data_2[31,2] = 0
$data_2[32,0] = 899705.9199999999
$data_2[32,1] = 1015751.04
$data_2[32,2] = 929467.008
$data_2[33,0] = 0
$data_2[33,1] = 0
$data_2[33,2] = 0
$data_2[34,0] = 0
$data_2[34,1] = 0
$data_2[34,2] = 0
$data_2[35,0] = 12661029.888
$data_2[35,1] = 12809053.184
$data_2[35,2] = 13076228.096
$data_2[36,0] = 10835383.296
$data_2[36,1] = 11083713.536
$data_2[36,2] = 11522839.552
$data_2[37,0] = 0
$data_2[37,1] = 0
$data_2[37,2] = 0
$data_2[38,0] = 282152
$data_2[38,1] = 228775.008
$data_2[38,2] = 223422
$data_2[39,0] = 397080
$data_2[39,1] = 470735.008
$data_2[39,2] = 526775.008
$data_2[40,0] = 0
$data_2[40,1] = 0
$data_2[40,2] = 0
$data_2[41,0] = 1146414.976
$data_2[41,1] = 1025828.992
$data_2[41,2] = 803190.976
$data_2[42,0] = 0
$data_2[42,1] = 0
$data_2[42,2] = 0
$data_2[43,0] = 0
$data_2[43,1] = 0
$data_2[43,2] = 0
$data_2[44,0] = 0
$data_2[44,1] = 0
$data_2[44,2] = 0
$data_2[45,0] = 4735294.976
$data_2[45,1] = 4889878.016
$data_2[45,2] = 5023527.936
$data_2[46,0] = 2142576
$data_2[46,1] = 2142576
$data_2[46,2] = 2142576
$data_2[47,0] = 1746089.984
$data_2[47,1] = 1714253.056
$data_2[47,2] = 1700630.016
$data_2[48,0] = 0
$data_2[48,1] = 0
$data_2[48,2] = 0
$data_2[49,0] = 865142.976
$data_2[49,1] = 1048160
$data_2[49,2] = 1189005.952
$data_2[50,0] = 0
$data_2[50,1] = 0
$data_2[50,2] = 0
$data_2[51,0] = 0
$data_2[51,1] = 0
$data_2[51,2] = 0
$data_2[52,0] = 0
$data_2[52,1] = 0
$data_2[52,2] = 0
$data_2[53,0] = -18514
$data_2[53,1] = -15111
$data_2[53,2] = -8684
$data_2[54,0] = 0
$data_2[54,1] = 0
$data_2[54,2] = 0
$ws.Range("V2:X56").Value = $data_2

# --- Numeric data rows 59:70 (V:X) ---
$data_59 = New-Object "object[,]" 12,3
$data_59[0,0] = 1452737.024
$data_59[0,1] = 1726110.976
$data_59[0,2] = 1883255.04
$data_59[1,0] = -749368.8959999999
$data_59[1,1] = -923137.024
$data_59[1,2] = -1023321.024
$data_59[2,0] = 703368.128
$data_59[2,1] = 802974.0159999999
$data_59[2,2] = 859934.0159999999
$data_59[3,0] = -74053
$data_59[3,1] = -58397
$data_59[3,2] = -68239
$data_59[4,0] = -101658.992
$data_59[4,1] = -81774
$data_59[4,2] = -91899
$data_59[5,0] = -31467
$data_59[5,1] = -25961
$data_59[5,2] = -109380
$data_59[6,0] = -13399
$data_59[6,1] = 4024
$data_59[6,2] = 0
$data_59[7,0] = 0
$data_59[7,1] = 0
$data_59[7,2] = -8281
$data_59[8,0] = 0
$data_59[8,1] = 0
$data_59[8,2] = 0
$data_59[9,0] = -395922.048
$data_59[9,1] = -395116
$data_59[9,2] = -407592
$data_59[10,0] = 48432.992
$data_59[10,1] = 61542
$data_59[10,2] = 54913
$data_59[11,0] = -444354.976
$data_59[11,1] = -456657.984
$data_59[11,2] = -462504.992
$ws.Range("V59:X70").Value = $data_59

# --- Numeric data rows 74:76 (V:X) ---
$data_74 = New-Object "object[,]" 3,3
$data_74[0,0] = 86867.992
$data_74[0,1] = 245750
$data_74[0,2] = 174543.008
$data_74[1,0] = 33573
$data_74[1,1] = -4971
$data_74[1,2] = -2871
$data_74[2,0] = 74991
$data_74[2,1] = -57762
$data_74[2,2] = -30826
$ws.Range("V74:X76").Value = $data_74

# --- Numeric data rows 79:80 (V:X) ---
$data_79 = New-Object "object[,]" 2,3
$data_79[0,0] = 0
$data_79[0,1] = 0
$data_79[0,2] = 0
$data_79[1,0] = 195431.968
$data_79[1,1] = 183016.992
$data_79[1,2] = 140846
$ws.Range("V79:X80").Value = $data_79

# --- Separator rows: blank (empty text) cells in V:X, styled like column U ---
$sepRows = @(57,58,71,72,73,77,78)
foreach ($r in $sepRows) {
    $ws.Range("V$r" + ":X$r").Formula = "'"
    $ws.Range("U$r").Copy()
    $ws.Range("V$r" + ":X$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
